# Apply the scraped-data refresh for 2023/poland_division-2_2023-2024.xlsx
#  - Rows 34-36 get their match data rotated (row34<-old row36, row35<-old row34, row36<-old row35)
#  - Rows 49-50 swap their match data
#  - Two brand-new match rows (134, 135) are appended at the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow {
    param(
        [int]$Row,
        [string]$Home,
        [int]$HomeGoals,
        [string]$Away,
        [int]$AwayGoals,
        [double]$HomeOpenOdds,
        [double]$HomeCloseOdds,
        [string]$HomeCloseDt,
        [double]$DrawOpenOdds,
        [double]$DrawCloseOdds,
        [string]$DrawCloseDt,
        [double]$AwayOpenOdds,
        [double]$AwayCloseOdds,
        [string]$AwayCloseDt,
        [string]$Url
    )

    $ws.Cells.Item($Row, 6).Value = $Home          # F home
    $ws.Cells.Item($Row, 7).Value = $HomeGoals      # G home_ft_gols
    $ws.Cells.Item($Row, 8).Value = $Away           # H away
    $ws.Cells.Item($Row, 9).Value = $AwayGoals      # I away_ft_gols

    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds  # J home_opening_odds
    # K (home_opening_data_hora) is unchanged

    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds # L home_closing_odds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseDt   # M home_closing_data_hora

    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds  # N draw_opening_odds
    # O (draw_opening_data_hora) is unchanged

    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds # P draw_closing_odds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseDt   # Q draw_closing_data_hora

    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds  # R away_opening_odds
    # S (away_opening_data_hora) is unchanged

    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds # T away_closing_odds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseDt   # U away_closing_data_hora

    $ws.Cells.Item($Row, 22).Value = $Url           # V url_partida
}

# --- Rows 34-36: rotate match data (K / opening-data-hora columns stay put) ---

Set-MatchRow 34 "Wisla Pulawy" 2 "Stomil Olsztyn" 1 `
    2.27 1.94 "13/08/2023 16:49" `
    3.23 3.47 "13/08/2023 16:53" `
    3.06 3.71 "13/08/2023 16:53" `
    "https://www.betexplorer.com/football/poland/division-2/wisla-pulawy-stomil-olsztyn/6Xey0QRr/"

Set-MatchRow 35 "Stezyca" 1 "Polonia Bytom" 1 `
    2.24 2.24 "13/08/2023 10:12" `
    3.15 3.17 "13/08/2023 15:01" `
    3.08 3.08 "13/08/2023 10:12" `
    "https://www.betexplorer.com/football/poland/division-2/stezyca-polonia-bytom/tjCd78Z8/"

Set-MatchRow 36 "Ol. Grudziadz" 1 "Lech Poznan II" 2 `
    1.9 1.68 "13/08/2023 15:37" `
    3.54 3.82 "13/08/2023 15:37" `
    3.58 4.66 "13/08/2023 15:37" `
    "https://www.betexplorer.com/football/poland/division-2/ol-grudziadz-lech-poznan/0fZ9q5s2/"

# --- Rows 49-50: swap match data (K / opening-data-hora columns stay put) ---

Set-MatchRow 49 "Pogon Siedlce" 3 "Sandecja Nowy S." 0 `
    2.28 2.77 "26/08/2023 16:56" `
    3.08 3.12 "26/08/2023 16:55" `
    3.09 2.56 "26/08/2023 16:56" `
    "https://www.betexplorer.com/football/poland/division-2/pogon-siedlce-sandecja-nowy-s/bT3PVrIQ/"

Set-MatchRow 50 "GKS Jastrzebie" 4 "Polonia Bytom" 2 `
    1.85 1.83 "26/08/2023 16:58" `
    3.42 3.6 "26/08/2023 16:58" `
    4.09 4.03 "26/08/2023 16:58" `
    "https://www.betexplorer.com/football/poland/division-2/gks-jastrzebie-polonia-bytom/464LWO2K/"

# --- New rows 134 & 135 ---
# Seed them from the last existing data row so they inherit the same styling
# (bold/bordered index cell, date-formatted match-date cell, ...), then overwrite content.

$lastRow = 133

$ws.Range("A$lastRow`:V$lastRow").Copy($ws.Range("A134:V134"))
$ws.Range("A$lastRow`:V$lastRow").Copy($ws.Range("A135:V135"))

$ws.Cells.Item(134, 1).Value = 133
$ws.Cells.Item(134, 2).Value = "poland"
$ws.Cells.Item(134, 3).Value = "division-2"
$ws.Cells.Item(134, 4).Value = "2023-2024"
$ws.Cells.Item(134, 5).Value = 45233.75
Set-MatchRow 134 "Polonia Bytom" 2 "Olimpia Elblag" 1 `
    2.48 2.12 "03/11/2023 17:24" `
    2.99 3.17 "03/11/2023 16:12" `
    2.65 3.52 "03/11/2023 17:24" `
    "https://www.betexplorer.com/football/poland/division-2/polonia-bytom-olimpia-elblag/z5wjEV4q/"
$ws.Cells.Item(134, 11).Value = "02/11/2023 06:12"  # K home_opening_data_hora
$ws.Cells.Item(134, 15).Value = "02/11/2023 06:12"  # O draw_opening_data_hora
$ws.Cells.Item(134, 19).Value = "02/11/2023 06:12"  # S away_opening_data_hora

$ws.Cells.Item(135, 1).Value = 134
$ws.Cells.Item(135, 2).Value = "poland"
$ws.Cells.Item(135, 3).Value = "division-2"
$ws.Cells.Item(135, 4).Value = "2023-2024"
$ws.Cells.Item(135, 5).Value = 45233.75
Set-MatchRow 135 "Wisla Pulawy" 2 "GKS Jastrzebie" 2 `
    1.83 2.16 "03/11/2023 17:25" `
    3.48 3.51 "03/11/2023 17:26" `
    3.7 3.1 "03/11/2023 17:25" `
    "https://www.betexplorer.com/football/poland/division-2/wisla-pulawy-gks-jastrzebie/QyyJBxqo/"
$ws.Cells.Item(135, 11).Value = "02/11/2023 06:12"  # K home_opening_data_hora
$ws.Cells.Item(135, 15).Value = "02/11/2023 06:12"  # O draw_opening_data_hora
$ws.Cells.Item(135, 19).Value = "02/11/2023 06:12"  # S away_opening_data_hora
